# The "Names" worksheet stores, for each plate well, the human-readable
# cargo-strand name (e.g. antiNelson_id1_h2_pos1). This edit:
#   1. Relabels the column-A header (A1) from the old regex description
#      to the simpler "name-side-position" label.
#   2. Strips the now-unused "_idN" segment out of every name in the
#      grid (antiNelson_id1_h2_pos1 -> antiNelson_h2_pos1, etc.)
#   3. Leaves the cursor/selection on cell B2 (instead of C10).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Names")

# Re-label the header cell that used to hold the naming-convention regex.
$ws.Range("A1").Value = "name-side-position"

# Walk every used cell (skipping column A, which holds row letters) and
# drop the "_id<digits>" token from any cargo name found there.
$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 2; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $current = $cell.Text
        if ($current -and $current -ne "") {
            $updated = $current -replace '_id\d+', ''
            if ($updated -ne $current) {
                $cell.Value = $updated
            }
        }
    }
}

# Move the visible selection to B2, matching the saved worksheet state.
$ws.Activate()
$ws.Range("B2").Select()
